$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1417
$ws.Range("I4").Value = 430.16666
$ws.Range("J4").Value = 3390.6667
$ws.Range("K4").Value = 430.16666
$ws.Range("L4").Value = 3390.6667
$ws.Range("M4").Value = -316.16666
$ws.Range("N4").Value = -3618.6667
$ws.Range("H40").Value = 1011.53845
$ws.Range("I40").Value = 1075
$ws.Range("K40").Value = 1075
$ws.Range("M40").Value = -900
$ws.Range("H137").Value = 1925038
$ws.Range("I137").Value = 3449755.8
$ws.Range("J137").Value = 2567.7827
$ws.Range("K137").Value = 10349267.4
$ws.Range("L137").Value = 7703.348100000001
$ws.Range("M137").Value = -10346717.4
$ws.Range("N137").Value = -12803.3481

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1240.0526
$ws.Range("I2").Value = 1105.1818
$ws.Range("J2").Value = 1425.5
$ws.Range("K2").Value = 1105.1818
$ws.Range("L2").Value = 1425.5
$ws.Range("M2").Value = -992.1818000000001
$ws.Range("N2").Value = -1651.5
$ws.Range("H32").Value = 20804.422
$ws.Range("I32").Value = 17843.89
$ws.Range("J32").Value = 30628
$ws.Range("K32").Value = 17843.89
$ws.Range("L32").Value = 30628
$ws.Range("M32").Value = -17556.89
$ws.Range("N32").Value = -31202
$ws.Range("H63").Value = 1867
$ws.Range("I63").Value = 1897.5
$ws.Range("J63").Value = 1806
$ws.Range("K63").Value = 1897.5
$ws.Range("L63").Value = 1806
$ws.Range("M63").Value = -1211.5
$ws.Range("N63").Value = -3178
$ws.Range("H66").Value = 1867
$ws.Range("I66").Value = 1897.5
$ws.Range("J66").Value = 1806
$ws.Range("K66").Value = 9487.5
$ws.Range("L66").Value = 9030
$ws.Range("M66").Value = -6055.5
$ws.Range("N66").Value = -15894
$ws.Range("H74").Value = 9582673
$ws.Range("I74").Value = 12383965
$ws.Range("K74").Value = 12383965
$ws.Range("M74").Value = -12383091
$ws.Range("H77").Value = 9582673
$ws.Range("I77").Value = 12383965
$ws.Range("K77").Value = 61919825
$ws.Range("M77").Value = -61915457
$ws.Range("H116").Value = 1240.0526
$ws.Range("I116").Value = 1105.1818
$ws.Range("J116").Value = 1425.5
$ws.Range("K116").Value = 1105.1818
$ws.Range("L116").Value = 1425.5
$ws.Range("M116").Value = 1188.8182
$ws.Range("N116").Value = -6013.5
$ws.Range("H122").Value = 8549324
$ws.Range("I122").Value = 1954.5
$ws.Range("J122").Value = 22225114
$ws.Range("K122").Value = 5863.5
$ws.Range("L122").Value = 66675342
$ws.Range("M122").Value = -3413.5
$ws.Range("N122").Value = -66680242

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1240.0526
$ws.Range("I3").Value = 1105.1818
$ws.Range("J3").Value = 1425.5
$ws.Range("K3").Value = 1105.1818
$ws.Range("L3").Value = 1425.5
$ws.Range("M3").Value = -991.1818000000001
$ws.Range("N3").Value = -1653.5
$ws.Range("H105").Value = 29414336
$ws.Range("I105").Value = 38464108
$ws.Range("J105").Value = 2575
$ws.Range("K105").Value = 38464108
$ws.Range("L105").Value = 2575
$ws.Range("M105").Value = -38462361
$ws.Range("N105").Value = -6069
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 29000
$ws.Range("J15").Value = 29000
$ws.Range("L15").Value = 29000
$ws.Range("N15").Value = -29340
$ws.Range("H31").Value = 2361.2683
$ws.Range("I31").Value = 1838.9
$ws.Range("J31").Value = 2858.762
$ws.Range("K31").Value = 1838.9
$ws.Range("L31").Value = 2858.762
$ws.Range("M31").Value = -1543.9
$ws.Range("N31").Value = -3448.762
$ws.Range("H34").Value = 2361.2683
$ws.Range("I34").Value = 1838.9
$ws.Range("J34").Value = 2858.762
$ws.Range("K34").Value = 1838.9
$ws.Range("L34").Value = 2858.762
$ws.Range("M34").Value = -1636.9
$ws.Range("N34").Value = -3262.762
$ws.Range("H105").Value = 607.2143
$ws.Range("I105").Value = 538.53845
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 538.53845
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 1208.46155
$ws.Range("N105").Value = -4994
$ws.Range("H132").Value = 17269.375
$ws.Range("I132").Value = 1356.4878
$ws.Range("J132").Value = 45635.824
$ws.Range("K132").Value = 4069.463400000001
$ws.Range("L132").Value = 136907.472
$ws.Range("M132").Value = -1539.463400000001
$ws.Range("N132").Value = -141967.472

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 461.06668
$ws.Range("I5").Value = 182.81818
$ws.Range("J5").Value = 1226.25
$ws.Range("K5").Value = 548.4545400000001
$ws.Range("L5").Value = 3678.75
$ws.Range("M5").Value = -436.4545400000001
$ws.Range("N5").Value = -3902.75
$ws.Range("H60").Value = 999
$ws.Range("J60").Value = 999
$ws.Range("L60").Value = 2997
$ws.Range("N60").Value = -3499
$ws.Range("H107").Value = 732.4706
$ws.Range("I107").Value = 904.0714
$ws.Range("J107").Value = 612.35
$ws.Range("K107").Value = 2712.2142
$ws.Range("L107").Value = 1837.05
$ws.Range("M107").Value = -792.2142000000003
$ws.Range("N107").Value = -5677.05
$ws.Range("H113").Value = 562.1539
$ws.Range("I113").Value = 510.5
$ws.Range("J113").Value = 734.3333
$ws.Range("K113").Value = 1531.5
$ws.Range("L113").Value = 2202.9999
$ws.Range("M113").Value = 638.5
$ws.Range("N113").Value = -6542.9999
$ws.Range("H131").Value = 1036.0344
$ws.Range("I131").Value = 399.875
$ws.Range("J131").Value = 1278.381
$ws.Range("K131").Value = 1199.625
$ws.Range("L131").Value = 3835.143
$ws.Range("M131").Value = 3840.375
$ws.Range("N131").Value = -13915.143
$ws.Range("H132").Value = 2829.4644
$ws.Range("I132").Value = 948
$ws.Range("J132").Value = 3342.5908
$ws.Range("K132").Value = 8532
$ws.Range("L132").Value = 30083.3172
$ws.Range("M132").Value = -6002
$ws.Range("N132").Value = -35143.3172
$ws.Range("H135").Value = 461.06668
$ws.Range("I135").Value = 182.81818
$ws.Range("J135").Value = 1226.25
$ws.Range("K135").Value = 1645.36362
$ws.Range("L135").Value = 11036.25
$ws.Range("M135").Value = 889.6363799999999
$ws.Range("N135").Value = -16106.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1920
$ws.Range("I82").Value = 1300
$ws.Range("K82").Value = 1300
$ws.Range("M82").Value = -939
$ws.Range("H85").Value = 1920
$ws.Range("I85").Value = 1300
$ws.Range("K85").Value = 1300
$ws.Range("M85").Value = -52
$ws.Range("H122").Value = 2983.3462
$ws.Range("I122").Value = 2712.9092
$ws.Range("K122").Value = 8138.7276
$ws.Range("M122").Value = -5688.7276
$ws.Range("H132").Value = 59634.777
$ws.Range("I132").Value = 4085.6667
$ws.Range("J132").Value = 170733
$ws.Range("K132").Value = 12257.0001
$ws.Range("L132").Value = 512199
$ws.Range("M132").Value = -9727.000100000001
$ws.Range("N132").Value = -517259

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 730.6
$ws.Range("I3").Value = 476.5
$ws.Range("J3").Value = 900
$ws.Range("K3").Value = 476.5
$ws.Range("L3").Value = 900
$ws.Range("M3").Value = -362.5
$ws.Range("N3").Value = -1128
$ws.Range("H122").Value = 1811.25
$ws.Range("I122").Value = 1435.625
$ws.Range("J122").Value = 2562.5
$ws.Range("K122").Value = 4306.875
$ws.Range("L122").Value = 7687.5
$ws.Range("M122").Value = -1856.875
$ws.Range("N122").Value = -12587.5
$ws.Range("H136").Value = 136291.73
$ws.Range("I136").Value = 168475.67
$ws.Range("K136").Value = 505427.01
$ws.Range("M136").Value = -502877.01
